$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.029050171712995
$ws.Range("D2").Value = 1.038593922272283
$ws.Range("E2").Value = 1.028939484054135
$ws.Range("F2").Value = 1.049556264529468
$ws.Range("I2").Value = 1.037227963336635
$ws.Range("J2").Value = 1.03419909610919
$ws.Range("K2").Value = 1.041381511715689
$ws.Range("L2").Value = 1.031754843321486
$ws.Range("M2").Value = 1.052312994001108
$ws.Range("N2").Value = 1.015412848883221
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029797076104069
$ws.Range("D3").Value = 1.039172919366744
$ws.Range("E3").Value = 1.029568752458215
$ws.Range("F3").Value = 1.050305086773869
$ws.Range("I3").Value = 1.037382267849378
$ws.Range("J3").Value = 1.034587905700218
$ws.Range("K3").Value = 1.041771164593182
$ws.Range("L3").Value = 1.032192634440216
$ws.Range("M3").Value = 1.052874249870005
$ws.Range("N3").Value = 1.015542531961086
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030281153918554
$ws.Range("D4").Value = 1.039548269894947
$ws.Range("E4").Value = 1.029976992166319
$ws.Range("F4").Value = 1.050790728854377
$ws.Range("I4").Value = 1.037481343308954
$ws.Range("J4").Value = 1.034839578222447
$ws.Range("K4").Value = 1.042023294375968
$ws.Range("L4").Value = 1.03247628446885
$ws.Range("M4").Value = 1.053237855903543
$ws.Range("N4").Value = 1.015626448408771
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030484845130413
$ws.Range("D5").Value = 1.039706233011509
$ws.Range("E5").Value = 1.030148868095853
$ws.Range("F5").Value = 1.050995154807279
$ws.Range("I5").Value = 1.037522809530752
$ws.Range("J5").Value = 1.034945400758891
$ws.Range("K5").Value = 1.042129287906403
$ws.Range("L5").Value = 1.032595617996749
$ws.Range("M5").Value = 1.05339081838611
$ws.Range("N5").Value = 1.015661727075722
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.03051905656531
$ws.Range("D6").Value = 1.039732765375769
$ws.Range("E6").Value = 1.030177741536076
$ws.Range("F6").Value = 1.051029494123021
$ws.Range("I6").Value = 1.037529761014444
$ws.Range("J6").Value = 1.034963169927163
$ws.Range("K6").Value = 1.042147084537715
$ws.Range("L6").Value = 1.032615659682547
$ws.Range("M6").Value = 1.053416507416816
$ws.Range("N6").Value = 1.015667650516334
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.03028387492649
$ws.Range("D7").Value = 1.039550379957051
$ws.Range("E7").Value = 1.029979287793462
$ws.Range("F7").Value = 1.050793459377269
$ws.Range("I7").Value = 1.037481898110948
$ws.Range("J7").Value = 1.034840992153797
$ws.Range("K7").Value = 1.042024710675695
$ws.Range("L7").Value = 1.032477878669108
$ws.Range("M7").Value = 1.053239899395286
$ws.Range("N7").Value = 1.015626919803898
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029302428685177
$ws.Range("D8").Value = 1.03878945077138
$ws.Range("E8").Value = 1.029151927348962
$ws.Range("F8").Value = 1.049809102547985
$ws.Range("I8").Value = 1.037280270201582
$ws.Range("J8").Value = 1.034330477362546
$ws.Range("K8").Value = 1.041513196221039
$ws.Range("L8").Value = 1.031902719388668
$ws.Range("M8").Value = 1.052502581662591
$ws.Range("N8").Value = 1.015456675010948
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027579063458512
$ws.Range("D9").Value = 1.037454050863971
$ws.Range("E9").Value = 1.02770222862042
$ws.Range("F9").Value = 1.04808309641818
$ws.Range("I9").Value = 1.036919114776184
$ws.Range("J9").Value = 1.033431609052808
$ws.Range("K9").Value = 1.040611892625273
$ws.Range("L9").Value = 1.030892114534558
$ws.Range("M9").Value = 1.05120675492646
$ws.Range("N9").Value = 1.015156724960627
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026434353815088
$ws.Range("D10").Value = 1.03656757088321
$ws.Range("E10").Value = 1.026741409159042
$ws.Range("F10").Value = 1.046938315194816
$ws.Range("I10").Value = 1.03667445119293
$ws.Range("J10").Value = 1.032832931270014
$ws.Range("K10").Value = 1.040011142818644
$ws.Range("L10").Value = 1.030220417047178
$ws.Range("M10").Value = 1.050345282347628
$ws.Range("N10").Value = 1.014956816147202
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025939703752315
$ws.Range("D11").Value = 1.03618463875436
$ws.Range("E11").Value = 1.026326728971373
$ws.Range("F11").Value = 1.046444038257575
$ws.Range("I11").Value = 1.036567595800246
$ws.Range("J11").Value = 1.032573849037747
$ws.Range("K11").Value = 1.039751057958042
$ws.Range("L11").Value = 1.029930067329833
$ws.Range("M11").Value = 1.049972850993387
$ws.Range("N11").Value = 1.01487027329958
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025756123531304
$ws.Range("D12").Value = 1.03604254093735
$ws.Range("E12").Value = 1.0261729048895
$ws.Range("F12").Value = 1.046260657644503
$ws.Range("I12").Value = 1.036527768372643
$ws.Range("J12").Value = 1.032477638174301
$ws.Range("K12").Value = 1.03965445891211
$ws.Range("L12").Value = 1.02982229515105
$ws.Range("M12").Value = 1.049834604274116
$ws.Range("N12").Value = 1.014838130784917
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02579549507182
$ws.Range("D13").Value = 1.03607301502271
$ws.Range("E13").Value = 1.026205891293811
$ws.Range("F13").Value = 1.046299983602747
$ws.Range("I13").Value = 1.036536317658471
$ws.Range("J13").Value = 1.032498274631076
$ws.Range("K13").Value = 1.039675179347348
$ws.Range("L13").Value = 1.029845409156436
$ws.Range("M13").Value = 1.049864254521691
$ws.Range("N13").Value = 1.014845025303713
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025924525783332
$ws.Range("D14").Value = 1.036172890031942
$ws.Range("E14").Value = 1.026314009597822
$ws.Range("F14").Value = 1.0464288755419
$ws.Range("I14").Value = 1.036564306434107
$ws.Range("J14").Value = 1.032565895723899
$ws.Range("K14").Value = 1.039743072877992
$ws.Range("L14").Value = 1.029921157275616
$ws.Range("M14").Value = 1.049961421611588
$ws.Range("N14").Value = 1.014867616318019
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026004046442796
$ws.Range("D15").Value = 1.036234444964626
$ws.Range("E15").Value = 1.026380652293918
$ws.Range("F15").Value = 1.046508318792729
$ws.Range("I15").Value = 1.036581533168939
$ws.Range("J15").Value = 1.032607562508242
$ws.Range("K15").Value = 1.039784905434015
$ws.Range("L15").Value = 1.029967838393091
$ws.Range("M15").Value = 1.050021301561304
$ws.Range("N15").Value = 1.014881535846667
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.026467203710203
$ws.Range("D16").Value = 1.036593004378441
$ws.Range("E16").Value = 1.026768958973723
$ws.Range("F16").Value = 1.046971148866692
$ws.Range("I16").Value = 1.036681523636047
$ws.Range("J16").Value = 1.032850128967294
$ws.Range("K16").Value = 1.040028404842451
$ws.Range("L16").Value = 1.030239697263204
$ws.Range("M16").Value = 1.050370012026968
$ws.Range("N16").Value = 1.014962560153689
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02675800371216
$ws.Range("D17").Value = 1.036818166881031
$ws.Range("E17").Value = 1.027012899374167
$ws.Range("F17").Value = 1.047261852034857
$ws.Range("I17").Value = 1.036744000761199
$ws.Range("J17").Value = 1.033002325430386
$ws.Range("K17").Value = 1.040181158348816
$ws.Range("L17").Value = 1.030410362040156
$ws.Range("M17").Value = 1.050588908738428
$ws.Range("N17").Value = 1.015013389997325
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026927720397136
$ws.Range("D18").Value = 1.036949588944035
$ws.Range("E18").Value = 1.027155316854686
$ws.Range("F18").Value = 1.047431551210265
$ws.Range("I18").Value = 1.036780354326165
$ws.Range("J18").Value = 1.03309111335415
$ws.Range("K18").Value = 1.040270261008984
$ws.Range("L18").Value = 1.030509956029367
$ws.Range("M18").Value = 1.050716644472521
$ws.Range("N18").Value = 1.015043039983733
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026985605953964
$ws.Range("D19").Value = 1.036994415433628
$ws.Range("E19").Value = 1.027203899719285
$ws.Range("F19").Value = 1.04748943739798
$ws.Range("I19").Value = 1.036792734953099
$ws.Range("J19").Value = 1.033121390124952
$ws.Range("K19").Value = 1.040300643368117
$ws.Range("L19").Value = 1.030543923119251
$ws.Range("M19").Value = 1.050760208663783
$ws.Range("N19").Value = 1.015053150152734
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02672679347622
$ws.Range("D20").Value = 1.036793999897743
$ws.Range("E20").Value = 1.026986713294717
$ws.Range("F20").Value = 1.047230648155064
$ws.Range("I20").Value = 1.03673730668287
$ws.Range("J20").Value = 1.032985994700891
$ws.Range("K20").Value = 1.040164768905997
$ws.Range("L20").Value = 1.030392046337056
$ws.Range("M20").Value = 1.05056541728467
$ws.Range("N20").Value = 1.015007936246514
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025886525163851
$ws.Range("D21").Value = 1.036143475443499
$ws.Range("E21").Value = 1.026282165731706
$ws.Range("F21").Value = 1.046390914100802
$ws.Range("I21").Value = 1.036556068200983
$ws.Range("J21").Value = 1.032545982328893
$ws.Range("K21").Value = 1.039723079691315
$ws.Range("L21").Value = 1.029898849214534
$ws.Range("M21").Value = 1.049932805801518
$ws.Range("N21").Value = 1.014860963730826
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025359111992212
$ws.Range("D22").Value = 1.035735277225577
$ws.Range("E22").Value = 1.025840385111055
$ws.Range("F22").Value = 1.045864189982682
$ws.Range("I22").Value = 1.03644132658247
$ws.Range("J22").Value = 1.032269467553673
$ws.Range("K22").Value = 1.039445419464667
$ws.Range("L22").Value = 1.029589200757912
$ws.Range("M22").Value = 1.04953558454465
$ws.Range("N22").Value = 1.014768575918751
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025638617986486
$ws.Range("D23").Value = 1.035951593134022
$ws.Range("E23").Value = 1.026074467261793
$ws.Range("F23").Value = 1.046143297124155
$ws.Range("I23").Value = 1.036502227878349
$ws.Range("J23").Value = 1.032416039676348
$ws.Range("K23").Value = 1.039592607427243
$ws.Range("L23").Value = 1.029753308644489
$ws.Range("M23").Value = 1.049746108479353
$ws.Range("N23").Value = 1.014817550417547
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026740895736003
$ws.Range("D24").Value = 1.036804919643811
$ws.Range("E24").Value = 1.026998545251851
$ws.Range("F24").Value = 1.047244747422406
$ws.Range("I24").Value = 1.036740331721621
$ws.Range("J24").Value = 1.032993373810604
$ws.Range("K24").Value = 1.04017217457719
$ws.Range("L24").Value = 1.030400322266452
$ws.Range("M24").Value = 1.050576031885634
$ws.Range("N24").Value = 1.015010400556263
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028023863131982
$ws.Range("D25").Value = 1.037798625051769
$ws.Range("E25").Value = 1.028076024222997
$ws.Range("F25").Value = 1.0485282823058
$ws.Range("I25").Value = 1.037013171814014
$ws.Range("J25").Value = 1.033663893911919
$ws.Range("K25").Value = 1.040844886443767
$ws.Range("L25").Value = 1.031153027275904
$ws.Range("M25").Value = 1.05154134024886
$ws.Range("N25").Value = 1.015234261153158
